# Testdata/Data.xlsx - "Login" sheet: add a Middle Name row (File Upload related
# fields being extended) and a Dob row; also nudge the view/selection to match
# what was left selected when the author saved.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Login")
$ws.Activate()

# --- Insert "Middle Name" / "Rajesh" as new row 7 (pushes Last Name..No. of
#     Dependents down by one row, matching the rest of the name fields) ---
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Middle Name"
$ws.Range("B7").Value = "Rajesh"

# --- Append "Dob" / "19/07/1981" as the new last row (14) ---
$ws.Range("A14").Value = "Dob"
$ws.Range("B14").Value = "19/07/1981"
$ws.Range("B14").NumberFormat = $ws.Range("B13").NumberFormat

# --- Leave the view scrolled down with G14 selected, like the saved file ---
try {
    $win = $excel.ActiveWindow
    $win.ScrollRow = 7
    $win.ScrollColumn = 1
} catch {
    # headless runtime may not expose window scroll position; non-fatal
}
$ws.Range("G14").Select()
